$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("230").Insert()

$ws.Range("A230").Value = 9
$ws.Range("B230").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C230").Value = "Metropolitana"
$ws.Range("D230").Value = 44588
$ws.Range("E230").Value = 13
$ws.Range("F230").Value = 100112052
$ws.Range("G230").Value = "Albahaca"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 160
$ws.Range("K230").Value = 3000
$ws.Range("L230").Value = 3500
$ws.Range("M230").Value = 3250
$ws.Range("N230").Value = "`$/docena de matas"
$ws.Range("O230").Value = "Región Metropolitana"
$ws.Range("P230").Value = 542
$ws.Range("Q230").Value = 6
$ws.Range("R230").Value = "Hortaliza"
